# Apply cell-level edits to Arkusz1 (sheet1) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Notatki (G) / Prototyp (H) column updates, row by row ---
$ws.Range("G2").Value = "Dane to nazwa, miejsce, opis, regulamin, zdjęcia, wyżywienie, zakwaterowanie, transport."
$ws.Range("H2").Value = "p3_dodawanie_edycja_imprezy"
$ws.Range("G3").Value = "Dane to nazwa, miejsce, opis, wyżywienie, zakwaterowanie, transport."
$ws.Range("H3").Value = "p2_wyszukanie_imprez"
$ws.Range("G4").Value = "Dane to nazwa, miejsce, opis, regulamin, zdjęcia, wyżywienie, zakwaterowanie, transport."
$ws.Range("H4").Value = "p3_dodawanie_edycja_imprezy"
$ws.Range("G5").Value = "Dane to nazwa, miejsce, opis, regulamin, zdjęcia, wyżywienie, zakwaterowanie, transport."
$ws.Range("H5").Value = "p4_usuwanie_imprezy"
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = "p7_dodawanie_edycja_skladnika"
$ws.Range("H7").Value = "p5_wyszukanie_skladnikow"
$ws.Range("H8").Value = "p7_dodawanie_edycja_skladnika"
$ws.Range("H9").Value = "p6_usuwanie_skladnika"
$ws.Range("G10").Value = ""
$ws.Range("H10").Value = "p9_dodawanie_edycja_katalogu"
$ws.Range("H11").Value = "p8_wyszukanie_katalogow"
$ws.Range("H12").Value = "p9_dodawanie_edycja_katalogu"
$ws.Range("H13").Value = "p10_usuwanie_katalogu"
$ws.Range("H14").Value = "p12_dodawanie_terminu"
$ws.Range("H15").Value = "p11_wyszukanie_terminow"
$ws.Range("H16").Value = "p13_edycja_terminu"
$ws.Range("H17").Value = "p14_usuwanie_terminow"
$ws.Range("H18").Value = "p21_dodawanie_edycja_cenniku"
$ws.Range("H19").Value = "p20_wyszukanie_cennikow"
$ws.Range("H20").Value = "p21_dodawanie_edycja_cenniku"
$ws.Range("H21").Value = "p22_usuwanie_cennikow"
$ws.Range("H22").Value = "p16_dodawanie_edycja_miasta_w_slowniku p17_dodawanie_edycja_hotelu_w_slowniku p18_dodawanie_edycja_panstwa_w_slowniku"
$ws.Range("H23").Value = "p15_wyszukanie_w_slowniku"
$ws.Range("H24").Value = "p16_dodawanie_edycja_miasta_w_slowniku p17_dodawanie_edycja_hotelu_w_slowniku p18_dodawanie_edycja_panstwa_w_slowniku"
$ws.Range("H25").Value = "p19_usuwanie_w_slowniku"
$ws.Range("H26").Value = "p2_wyszukanie_imprez"
$ws.Range("H27").Value = "p1_homepage p23_wyszukanie_imprez_przez_klienta p24_podglad_imprezy_klient"
$ws.Range("H28").Value = "p23_wyszukanie_imprez_przez_klienta "
$ws.Range("H29").Value = "p26_cennik_imprezy"
$ws.Range("H30").Value = "p25_rezerwacja_imprezy"
$ws.Range("H31").Value = "p27_platnosc_po_rezerwacji"
$ws.Range("H32").Value = "p28_rezygnacja"

# --- Row heights for rows whose wrapped content grew taller ---
$ws.Rows.Item(22).RowHeight = 90
$ws.Rows.Item(24).RowHeight = 90
$ws.Rows.Item(27).RowHeight = 60

# --- Selection moves from I7 to G6 ---
$ws.Range("G6").Select()
